# Remove the "Header Location" column (column G) from the link-check report.
# Since no successful pages are reported, the header-location value is never
# populated, so the whole column is dropped and everything to its right
# shifts one column to the left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Deleting the entire column removes the cells, shifts H..K left to G..J,
# and updates the shared-string usage / dimension accordingly.
$ws.Columns.Item(7).Delete()

# Re-apply the AutoFilter so its range shrinks from A1:H11 to A1:G11 along
# with the sheet's used range.
$ws.AutoFilterMode = $false
$ws.Range("A1:G11").AutoFilter()

# The workbook-level _FilterDatabase defined name also needs to track the
# new, narrower range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$G`$11"
    }
}

# Move the active selection off the (now gone) old C3 cell onto H1, matching
# the post-edit cursor position recorded by Excel.
$ws.Range("H1").Select()
